{"js": "// Office.js (Word JavaScript API) script.\n//\n// Change described by the diff:\n//   - The paragraph ending in \"...websockets.\" (which carries the\n//     _GoBack bookmark) is split right after the final period.\n//   - A new paragraph \"Torsdag:\" is inserted.\n//   - A new paragraph with the day's status text is inserted; the\n//     _GoBack bookmark moves to the end of this new paragraph.\n//   - The old (already existing) \"Torsdag:\" paragraph is left as-is.\n//   - The empty paragraph that used to sit between the old \"Torsdag:\"\n//     and \"Fredag:\" paragraphs is removed.\n\n// --- Step 1: locate the end of \"websockets.\" (end of the paragraph that\n// currently holds the _GoBack bookmark). ---------------------------------\nconst found = context.document.body.search(\"websockets.\", { matchCase: true });\nfound.load(\"text\");\nawait context.sync();\n\nconst afterWebsockets = found.items[0].getRange(\"After\");\n\n// --- Step 2: insert \"Torsdag:\" then the new status paragraph -------------\nconst torsdagPara = afterWebsockets.insertParagraph(\"Torsdag:\", \"After\");\nawait context.sync();\n\nconst statusPara = torsdagPara\n  .getRange(\"End\")\n  .insertParagraph(\n    \"Lavede hjemmesiden og lokalt program s\u00e5 det eneste der mangler er at koble SQL server p\u00e5.\",\n    \"After\"\n  );\nawait context.sync();\n\n// --- Step 3: move the _GoBack bookmark to the end of the new status\n// paragraph (delete the old one, re-insert at the new location). ---------\nconst hasGoBack = context.document.bookmarks.exists(\"_GoBack\");\nawait context.sync();\nif (hasGoBack.value) {\n  context.document.deleteBookmark(\"_GoBack\");\n}\nstatusPara.getRange(\"End\").insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- Step 4: remove the old empty paragraph that used to sit between the\n// (original) \"Torsdag:\" paragraph and \"Fredag:\". --------------------------\nconst fredagFound = context.document.body.search(\"Fredag:\", { matchCase: true });\nfredagFound.load(\"text\");\nawait context.sync();\n\nconst fredagPara = fredagFound.items[0].paragraphs.getFirst();\nfredagPara.load(\"text\");\nawait context.sync();\n\nconst prevPara = fredagPara.getPrevious();\nprevPara.load(\"text\");\nawait context.sync();\n\nif (prevPara.text.trim() === \"\") {\n  prevPara.delete();\n  await context.sync();\n}\n", "ps1": "# Word COM interop script.\n#\n# Change described by the diff:\n#   - The paragraph ending in \"...websockets.\" (which carries the\n#     _GoBack bookmark) is split right after the final period.\n#   - A new paragraph \"Torsdag:\" is inserted.\n#   - A new paragraph with the day's status text is inserted; the\n#     _GoBack bookmark moves to the end of this new paragraph.\n#   - The old (already existing) \"Torsdag:\" paragraph is left as-is.\n#   - The empty paragraph that used to sit between the old \"Torsdag:\"\n#     and \"Fredag:\" paragraphs is removed.\n\n$d = $word.ActiveDocument\n\n$wdParagraph    = 4\n$aa             = [char]0x00E5   # '\u00e5'\n\n# --- Step 1: locate the paragraph that ends in \"websockets.\" -------------\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"websockets.\"\n$null = $rng.Find.Execute()\n$null = $rng.Expand($wdParagraph)\n$insertPos = $rng.End - 1   # position right before the paragraph mark\n\n# --- Step 2: drop the existing _GoBack bookmark ---------------------------\n# (it gets recreated below, anchored at the end of the new status text)\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# --- Step 3: insert \"Torsdag:\" + the new status paragraph ----------------\n$insertRng = $d.Range($insertPos, $insertPos)\n$newText = \"`rTorsdag:`rLavede hjemmesiden og lokalt program s\" + $aa + \" det eneste der mangler er at koble SQL server p\" + $aa + \".\"\n# Insert a trailing placeholder character too: adding a zero-length\n# bookmark exactly at an end-of-paragraph offset is unreliable, so we\n# anchor the bookmark one character early (non end-of-paragraph) and\n# then delete the placeholder.\n$insertRng.InsertAfter($newText + \"X\")\n\n# --- Step 4: recreate _GoBack at the end of the new status paragraph -----\n$bmPos = $insertRng.End - 1\n$bmRng = $d.Range($bmPos, $bmPos)\n$d.Bookmarks.Add(\"_GoBack\", $bmRng)\n\n# --- Step 5: remove the placeholder character -----------------------------\n$d.Range($bmPos, $bmPos + 1).Delete()\n\n# --- Step 6: remove the old empty paragraph before \"Fredag:\" -------------\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Text = \"Fredag:\"\n$null = $rng2.Find.Execute()\n$null = $rng2.Expand($wdParagraph)\n\n$prevPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Start -eq $rng2.Start) {\n        $prevPara = $d.Paragraphs.Item($i - 1)\n        break\n    }\n}\nif ($prevPara -ne $null -and $prevPara.Range.Text.Trim() -eq \"\") {\n    $prevPara.Range.Delete()\n}\n"}
